$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 3675

$ws.Cells.Item(73, 8).Value = 3675

$ws.Cells.Item(80, 8).Value = 19480.908
$ws.Cells.Item(80, 10).Value = 18049.666
$ws.Cells.Item(80, 12).Value = 54148.99800000001
$ws.Cells.Item(80, 14).Value = -56144.99800000001

$ws.Cells.Item(83, 8).Value = 19480.908
$ws.Cells.Item(83, 10).Value = 18049.666
$ws.Cells.Item(83, 12).Value = 162446.994
$ws.Cells.Item(83, 14).Value = -172430.994

$ws.Cells.Item(100, 8).Value = 1359.3077
$ws.Cells.Item(100, 9).Value = 1305.4166
$ws.Cells.Item(100, 11).Value = 1305.4166
$ws.Cells.Item(100, 13).Value = -764.4166

$ws.Cells.Item(101, 8).Value = 4117.2383
$ws.Cells.Item(101, 9).Value = 1378.091
$ws.Cells.Item(101, 10).Value = 7130.3
$ws.Cells.Item(101, 11).Value = 4134.272999999999
$ws.Cells.Item(101, 12).Value = 21390.9
$ws.Cells.Item(101, 13).Value = -2512.272999999999
$ws.Cells.Item(101, 14).Value = -24634.9

$ws.Cells.Item(106, 8).Value = 4381
$ws.Cells.Item(106, 9).Value = 3999
$ws.Cells.Item(106, 11).Value = 3999
$ws.Cells.Item(106, 13).Value = -3368

$ws.Cells.Item(132, 8).Value = 1580.2727
$ws.Cells.Item(132, 9).Value = 1172.3513
$ws.Cells.Item(132, 11).Value = 3517.0539
$ws.Cells.Item(132, 13).Value = -987.0538999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(11, 8).Value = 1500
$ws.Cells.Item(11, 10).Value = 1500
$ws.Cells.Item(11, 12).Value = 1500
$ws.Cells.Item(11, 14).Value = -1788

$ws.Cells.Item(110, 8).Value = 3660.4167
$ws.Cells.Item(110, 9).Value = 2749
$ws.Cells.Item(110, 11).Value = 2749
$ws.Cells.Item(110, 13).Value = -704

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(52, 8).Value = 45000
$ws.Cells.Item(52, 10).Value = 45000
$ws.Cells.Item(52, 12).Value = 45000
$ws.Cells.Item(52, 14).Value = -45526

$ws.Cells.Item(121, 8).Value = 45000
$ws.Cells.Item(121, 10).Value = 45000
$ws.Cells.Item(121, 12).Value = 45000
$ws.Cells.Item(121, 14).Value = -48494

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 828
$ws.Cells.Item(6, 9).Value = 1293.4
$ws.Cells.Item(6, 10).Value = 246.25
$ws.Cells.Item(6, 11).Value = 1293.4
$ws.Cells.Item(6, 12).Value = 246.25
$ws.Cells.Item(6, 13).Value = -1180.4
$ws.Cells.Item(6, 14).Value = -472.25

$ws.Cells.Item(41, 8).Value = 22821.428
$ws.Cells.Item(41, 10).Value = 20558.334
$ws.Cells.Item(41, 12).Value = 20558.334
$ws.Cells.Item(41, 14).Value = -21414.334

$ws.Cells.Item(50, 8).Value = 29726.666
$ws.Cells.Item(50, 10).Value = 29726.666
$ws.Cells.Item(50, 12).Value = 29726.666
$ws.Cells.Item(50, 14).Value = -30976.666

$ws.Cells.Item(51, 8).Value = 27200
$ws.Cells.Item(51, 10).Value = 27200
$ws.Cells.Item(51, 12).Value = 27200
$ws.Cells.Item(51, 14).Value = -28672

$ws.Cells.Item(59, 8).Value = 15413.833
$ws.Cells.Item(59, 10).Value = 10496.6
$ws.Cells.Item(59, 12).Value = 10496.6
$ws.Cells.Item(59, 14).Value = -12786.6

$ws.Cells.Item(60, 8).Value = 19827.3
$ws.Cells.Item(60, 10).Value = 22154.6
$ws.Cells.Item(60, 12).Value = 22154.6
$ws.Cells.Item(60, 14).Value = -23176.6

$ws.Cells.Item(61, 8).Value = 27200
$ws.Cells.Item(61, 10).Value = 27200
$ws.Cells.Item(61, 12).Value = 27200
$ws.Cells.Item(61, 14).Value = -27896

$ws.Cells.Item(74, 8).Value = 39449.668
$ws.Cells.Item(74, 10).Value = 39449.668
$ws.Cells.Item(74, 12).Value = 39449.668
$ws.Cells.Item(74, 14).Value = -41197.668

$ws.Cells.Item(77, 8).Value = 39449.668
$ws.Cells.Item(77, 10).Value = 39449.668
$ws.Cells.Item(77, 12).Value = 118349.004
$ws.Cells.Item(77, 14).Value = -127085.004

$ws.Cells.Item(97, 8).Value = 25452.715
$ws.Cells.Item(97, 9).Value = 16000
$ws.Cells.Item(97, 11).Value = 16000
$ws.Cells.Item(97, 13).Value = -15009

$ws.Cells.Item(132, 8).Value = 5743.6875
$ws.Cells.Item(132, 9).Value = 4428.25
$ws.Cells.Item(132, 11).Value = 13284.75
$ws.Cells.Item(132, 13).Value = -10754.75

$ws.Cells.Item(140, 8).Value = 100486.11
$ws.Cells.Item(140, 10).Value = 100486.11
$ws.Cells.Item(140, 12).Value = 100486.11
$ws.Cells.Item(140, 14).Value = -110846.11

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 82105464
$ws.Cells.Item(4, 9).Value = 95739950
$ws.Cells.Item(4, 11).Value = 287219850
$ws.Cells.Item(4, 13).Value = -287219738

$ws.Cells.Item(107, 8).Value = 333.57144
$ws.Cells.Item(107, 9).Value = 296.66666
$ws.Cells.Item(107, 11).Value = 889.9999799999999
$ws.Cells.Item(107, 13).Value = 1030.00002

$ws.Cells.Item(140, 8).Value = 2602.375
$ws.Cells.Item(140, 9).Value = 2602.375
$ws.Cells.Item(140, 11).Value = 7807.125
$ws.Cells.Item(140, 13).Value = -2627.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value = 16562
$ws.Cells.Item(46, 10).Value = 17976.25
$ws.Cells.Item(46, 12).Value = 17976.25
$ws.Cells.Item(46, 14).Value = -18288.25

$ws.Cells.Item(126, 8).Value = 8071.143
$ws.Cells.Item(126, 9).Value = 5500
$ws.Cells.Item(126, 11).Value = 16500
$ws.Cells.Item(126, 13).Value = -14030

$ws.Cells.Item(132, 8).Value = 7240.6875
$ws.Cells.Item(132, 10).Value = 7969.25
$ws.Cells.Item(132, 12).Value = 23907.75
$ws.Cells.Item(132, 14).Value = -28967.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3186.7368
$ws.Cells.Item(7, 9).Value = 3097.125
$ws.Cells.Item(7, 10).Value = 3664.6667
$ws.Cells.Item(7, 11).Value = 3097.125
$ws.Cells.Item(7, 12).Value = 3664.6667
$ws.Cells.Item(7, 13).Value = -2985.125
$ws.Cells.Item(7, 14).Value = -3888.6667

$ws.Cells.Item(22, 8).Value = 5799.3335
$ws.Cells.Item(22, 10).Value = 5749.5
$ws.Cells.Item(22, 12).Value = 5749.5
$ws.Cells.Item(22, 14).Value = -6339.5

$ws.Cells.Item(27, 8).Value = 5799.3335
$ws.Cells.Item(27, 10).Value = 5749.5
$ws.Cells.Item(27, 12).Value = 5749.5
$ws.Cells.Item(27, 14).Value = -5963.5

$ws.Cells.Item(68, 8).Value = 5486.4
$ws.Cells.Item(68, 9).Value = 4108
$ws.Cells.Item(68, 10).Value = 11000
$ws.Cells.Item(68, 11).Value = 4108
$ws.Cells.Item(68, 12).Value = 11000
$ws.Cells.Item(68, 13).Value = -3359
$ws.Cells.Item(68, 14).Value = -12498

$ws.Cells.Item(71, 8).Value = 5486.4
$ws.Cells.Item(71, 9).Value = 4108
$ws.Cells.Item(71, 10).Value = 11000
$ws.Cells.Item(71, 11).Value = 20540
$ws.Cells.Item(71, 12).Value = 55000
$ws.Cells.Item(71, 13).Value = -16796
$ws.Cells.Item(71, 14).Value = -62488

$ws.Cells.Item(82, 8).Value = 1043.75
$ws.Cells.Item(82, 10).Value = 763.75
$ws.Cells.Item(82, 12).Value = 763.75
$ws.Cells.Item(82, 14).Value = -1485.75

$ws.Cells.Item(85, 8).Value = 1043.75
$ws.Cells.Item(85, 10).Value = 763.75
$ws.Cells.Item(85, 12).Value = 763.75
$ws.Cells.Item(85, 14).Value = -3259.75

$ws.Cells.Item(126, 8).Value = 3186.7368
$ws.Cells.Item(126, 9).Value = 3097.125
$ws.Cells.Item(126, 10).Value = 3664.6667
$ws.Cells.Item(126, 11).Value = 9291.375
$ws.Cells.Item(126, 12).Value = 10994.0001
$ws.Cells.Item(126, 13).Value = -6821.375
$ws.Cells.Item(126, 14).Value = -15934.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(51, 8).Value = 30077
$ws.Cells.Item(51, 10).Value = 30077
$ws.Cells.Item(51, 12).Value = 30077
$ws.Cells.Item(51, 14).Value = -31097

$ws.Cells.Item(54, 8).Value = 79460
$ws.Cells.Item(54, 10).Value = 99113.336
$ws.Cells.Item(54, 12).Value = 99113.336
$ws.Cells.Item(54, 14).Value = -100153.336

$ws.Cells.Item(81, 8).Value = 1428.5555
$ws.Cells.Item(81, 10).Value = 1000
$ws.Cells.Item(81, 12).Value = 2000
$ws.Cells.Item(81, 14).Value = -4122

$ws.Cells.Item(84, 8).Value = 1428.5555
$ws.Cells.Item(84, 10).Value = 1000
$ws.Cells.Item(84, 12).Value = 10000
$ws.Cells.Item(84, 14).Value = -20608

$ws.Cells.Item(93, 8).Value = 27475
$ws.Cells.Item(93, 9).Value = 27475
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 27475
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = -24979
$ws.Cells.Item(93, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 4248.0605
$ws.Cells.Item(132, 9).Value = 3831.698
$ws.Cells.Item(132, 10).Value = 5945.5386
$ws.Cells.Item(132, 11).Value = 11495.094
$ws.Cells.Item(132, 12).Value = 17836.6158
$ws.Cells.Item(132, 13).Value = -8965.093999999999
$ws.Cells.Item(132, 14).Value = -22896.6158

$ws.Cells.Item(136, 8).Value = 6929.7144
$ws.Cells.Item(136, 9).Value = 5584.6665
$ws.Cells.Item(136, 11).Value = 16753.9995
$ws.Cells.Item(136, 13).Value = -14203.9995
